# Auto-update predictions and index for 2025-10-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column E hold "X of Y Tips" text - not number-like, so plain
# assignment keeps them stored as text (matches the original inline-string
# representation).
$ws.Range("E2").Value = "85 of 90 Tips"
$ws.Range("E3").Value = "63 of 68 Tips"
$ws.Range("E4").Value = "54 of 67 Tips"
$ws.Range("E5").Value = "48 of 54 Tips"
$ws.Range("E6").Value = "45 of 55 Tips"
$ws.Range("E7").Value = "40 of 48 Tips"
$ws.Range("E8").Value = "33 of 43 Tips"
$ws.Range("E11").Value = "17 of 17 Tips"

# Cells in columns F/G hold numeric-looking text ("81", "1.28", ...) that
# was stored as text in the source workbook (inline string, General format).
# Assigning a numeric-looking string directly would auto-convert the cell
# to a real number, so force text formatting first, then restore the
# original (General/default) style by copying it from an untouched
# neighboring cell - this keeps the value textual without leaving any
# lingering formatting difference behind.
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1.28"
$ws.Range("G3").Style = $ws.Range("G2").Style

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "81"
$ws.Range("F4").Style = $ws.Range("A4").Style

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "82"
$ws.Range("F6").Style = $ws.Range("A6").Style

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "77"
$ws.Range("F8").Style = $ws.Range("A8").Style
